$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: relabel the "item group" header columns and add a 3rd item group ---
# Group 1 (ItemName2/Price2 -> ItemName1/Price1, plus new Discount1/Final Price1)
$ws1.Cells.Item(1, 8).Value  = "ItemName1"     # H1
$ws1.Cells.Item(1, 9).Value  = "Price1"        # I1

# New 3rd item group, columns R:U (reuses ItemName3/Price3 text, adds new Discout3/FinalPrice3)
$ws1.Cells.Item(1, 18).Value = "ItemName3"     # R1
$ws1.Cells.Item(1, 19).Value = "Price3"        # S1
$ws1.Cells.Item(1, 20).Value = "Discout3"      # T1
$ws1.Cells.Item(1, 21).Value = "FinalPrice3"   # U1

# Group 2 (relabel Discount/Final Price as Discount2/Final Price2)
$ws1.Cells.Item(1, 16).Value = "Final Price2"  # P1
$ws1.Cells.Item(1, 15).Value = "Discount2"     # O1

# Group 1 Discount/Final Price columns
$ws1.Cells.Item(1, 10).Value = "Discount1"     # J1
$ws1.Cells.Item(1, 11).Value = "Final Price1"  # K1

# Group 2 (ItemName3/Price3 -> ItemName2/Price2)
$ws1.Cells.Item(1, 13).Value = "ItemName2"     # M1
$ws1.Cells.Item(1, 14).Value = "Price2"        # N1

# --- Sheet2: same header relabeling ---
$ws2.Cells.Item(1, 4).Value  = "ItemName1"     # D1
$ws2.Cells.Item(1, 5).Value  = "Price1"        # E1

$ws2.Cells.Item(1, 10).Value = "ItemName3"     # J1
$ws2.Cells.Item(1, 11).Value = "Price3"        # K1

$ws2.Cells.Item(1, 7).Value  = "ItemName2"     # G1
$ws2.Cells.Item(1, 8).Value  = "Price2"        # H1

# --- Restore selections on each sheet, leaving Sheet1 as the active tab ---
[void]$ws2.Range("D6").Select()
[void]$ws1.Range("J11").Select()
